# TestDataMappingSheet.xlsx — "Add files via upload" re-edit.
#
# The rows for the "Folio.xlsx" block (A43:D47) got re-ordered: the last
# row's ScreenName/TestDataSheetName ("FolioERR") moved up to the top of
# the block (row 43) and everything else shifted down by one row. D
# (TotalIterations) stays 1 throughout, so only columns A and C change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestDataMappingSheet_SD")

# New column A / C values for rows 43-47 (shared strings already present
# in the workbook: Documents, Assignments, Folioaddendums, FolioDocuments,
# FolioERR — COM will reuse/append to the shared string table as needed).
$newValues = @{
    43 = "FolioERR"
    44 = "Documents"
    45 = "Assignments"
    46 = "Folioaddendums"
    47 = "FolioDocuments"
}

foreach ($row in $newValues.Keys) {
    $val = $newValues[$row]
    $ws.Cells.Item($row, 1).Value2 = $val   # column A
    $ws.Cells.Item($row, 3).Value2 = $val   # column C
}

# Match the author's final selection / scroll state: active cell A43,
# selected range A43:D43 (the first row of the re-ordered block).
$ws.Range("A43:D43").Select()
